$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168259620666504
$ws.Range("B1").Value = 2.280840158462524
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.326796054840088
$ws.Range("E1").Value = 1.232401013374329
